# Auto-generated edit script: updates computed price/profit columns (H-N)
# on multiple job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to refreshed
# market-board values, per the scheduled runner sync.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1608.4667
$ws.Range("I2").Value = 804.8570999999999
$ws.Range("J2").Value = 2311.625
$ws.Range("K2").Value = 804.8570999999999
$ws.Range("L2").Value = 2311.625
$ws.Range("M2").Value = -691.8570999999999
$ws.Range("N2").Value = -2537.625
$ws.Range("H40").Value = 4326.6665
$ws.Range("I40").Value = 2993
$ws.Range("K40").Value = 2993
$ws.Range("M40").Value = -2818
$ws.Range("H64").Value = 4714.2856
$ws.Range("I64").Value = 3400
$ws.Range("K64").Value = 3400
$ws.Range("M64").Value = -3152
$ws.Range("H67").Value = 4714.2856
$ws.Range("I67").Value = 3400
$ws.Range("K67").Value = 3400
$ws.Range("M67").Value = -2542
$ws.Range("H74").Value = 3858.4
$ws.Range("I74").Value = 3264.3333
$ws.Range("K74").Value = 3264.3333
$ws.Range("M74").Value = -2328.3333
$ws.Range("H77").Value = 3858.4
$ws.Range("I77").Value = 3264.3333
$ws.Range("K77").Value = 16321.6665
$ws.Range("M77").Value = -11641.6665
$ws.Range("H107").Value = 36969.953
$ws.Range("I107").Value = 928.8
$ws.Range("J107").Value = 127072.836
$ws.Range("K107").Value = 928.8
$ws.Range("L107").Value = 127072.836
$ws.Range("M107").Value = 991.2
$ws.Range("N107").Value = -130912.836
$ws.Range("H114").Value = 41472
$ws.Range("J114").Value = 41472
$ws.Range("L114").Value = 41472
$ws.Range("N114").Value = -50150
$ws.Range("H116").Value = 4552.778
$ws.Range("I116").Value = 4244
$ws.Range("K116").Value = 4244
$ws.Range("M116").Value = -802
$ws.Range("H138").Value = 2506.9473
$ws.Range("I138").Value = 3381.8333
$ws.Range("K138").Value = 10145.4999
$ws.Range("M138").Value = -5005.499899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1272.5652
$ws.Range("I2").Value = 1093.7
$ws.Range("K2").Value = 1093.7
$ws.Range("M2").Value = -980.7
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H74").Value = 1812.7742
$ws.Range("I74").Value = 1704
$ws.Range("K74").Value = 1704
$ws.Range("M74").Value = -830
$ws.Range("H77").Value = 1812.7742
$ws.Range("I77").Value = 1704
$ws.Range("K77").Value = 8520
$ws.Range("M77").Value = -4152
$ws.Range("H94").Value = 30329.666
$ws.Range("J94").Value = 30329.666
$ws.Range("L94").Value = 30329.666
$ws.Range("N94").Value = -32131.666
$ws.Range("H116").Value = 1272.5652
$ws.Range("I116").Value = 1093.7
$ws.Range("K116").Value = 1093.7
$ws.Range("M116").Value = 1200.3
$ws.Range("H132").Value = 3393.9285
$ws.Range("I132").Value = 1390.5555
$ws.Range("K132").Value = 4171.666499999999
$ws.Range("M132").Value = -1641.666499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1272.5652
$ws.Range("I3").Value = 1093.7
$ws.Range("K3").Value = 1093.7
$ws.Range("M3").Value = -979.7
$ws.Range("H105").Value = 13152.25
$ws.Range("I105").Value = 13152.25
$ws.Range("K105").Value = 13152.25
$ws.Range("M105").Value = -11405.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1175.25
$ws.Range("I16").Value = 1233.5714
$ws.Range("K16").Value = 1233.5714
$ws.Range("M16").Value = -946.5714
$ws.Range("H31").Value = 32024.777
$ws.Range("I31").Value = 28481.236
$ws.Range("J31").Value = 51261.145
$ws.Range("K31").Value = 28481.236
$ws.Range("L31").Value = 51261.145
$ws.Range("M31").Value = -28186.236
$ws.Range("N31").Value = -51851.145
$ws.Range("H34").Value = 32024.777
$ws.Range("I34").Value = 28481.236
$ws.Range("J34").Value = 51261.145
$ws.Range("K34").Value = 28481.236
$ws.Range("L34").Value = 51261.145
$ws.Range("M34").Value = -28279.236
$ws.Range("N34").Value = -51665.145
$ws.Range("H52").Value = 49974.5
$ws.Range("I52").Value = 49974.5
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 49974.5
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -49680.5
$ws.Range("N52").ClearContents()
$ws.Range("H113").Value = 1175.25
$ws.Range("I113").Value = 1233.5714
$ws.Range("K113").Value = 1233.5714
$ws.Range("M113").Value = 936.4286
$ws.Range("H132").Value = 2324.2942
$ws.Range("I132").Value = 2443.0435
$ws.Range("K132").Value = 7329.130500000001
$ws.Range("M132").Value = -4799.130500000001
$ws.Range("H134").Value = 2043.1765
$ws.Range("I134").Value = 1649
$ws.Range("K134").Value = 4947
$ws.Range("M134").Value = -2412

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 166.5
$ws.Range("I7").Value = 111
$ws.Range("J7").Value = 222
$ws.Range("K7").Value = 333
$ws.Range("L7").Value = 666
$ws.Range("M7").Value = -221
$ws.Range("N7").Value = -890
$ws.Range("H33").Value = 599.7778
$ws.Range("I33").Value = 199.33333
$ws.Range("J33").Value = 800
$ws.Range("K33").Value = 1195.99998
$ws.Range("L33").Value = 4800
$ws.Range("M33").Value = -912.9999800000001
$ws.Range("N33").Value = -5366
$ws.Range("H39").Value = 108568.6
$ws.Range("I39").Value = 999999
$ws.Range("J39").Value = 9520.777
$ws.Range("K39").Value = 2999997
$ws.Range("L39").Value = 28562.331
$ws.Range("M39").Value = -2999703
$ws.Range("N39").Value = -29150.331
$ws.Range("H56").Value = 8275.84
$ws.Range("I56").Value = 8275.84
$ws.Range("K56").Value = 8275.84
$ws.Range("M56").Value = -7745.84

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4819.75
$ws.Range("I70").Value = 4562.706
$ws.Range("J70").Value = 5217
$ws.Range("K70").Value = 4562.706
$ws.Range("L70").Value = 5217
$ws.Range("M70").Value = -4292.706
$ws.Range("N70").Value = -5757
$ws.Range("H73").Value = 4819.75
$ws.Range("I73").Value = 4562.706
$ws.Range("J73").Value = 5217
$ws.Range("K73").Value = 4562.706
$ws.Range("L73").Value = 5217
$ws.Range("M73").Value = -3626.706
$ws.Range("N73").Value = -7089
$ws.Range("H102").Value = 1621
$ws.Range("I102").Value = 1457.6522
$ws.Range("K102").Value = 1457.6522
$ws.Range("M102").Value = 164.3478
$ws.Range("H132").Value = 4166.3335
$ws.Range("I132").Value = 4166.3335
$ws.Range("K132").Value = 12499.0005
$ws.Range("M132").Value = -9969.000499999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 220149.6
$ws.Range("J16").Value = 850
$ws.Range("L16").Value = 850
$ws.Range("N16").Value = -1190
$ws.Range("H46").Value = 23129.285
$ws.Range("I46").Value = 71736.836
$ws.Range("K46").Value = 71736.836
$ws.Range("M46").Value = -71548.836
$ws.Range("H61").Value = 107746
$ws.Range("I61").Value = 106071.25
$ws.Range("J61").Value = 114445
$ws.Range("K61").Value = 106071.25
$ws.Range("L61").Value = 114445
$ws.Range("M61").Value = -105869.25
$ws.Range("N61").Value = -114849
$ws.Range("H113").Value = 107746
$ws.Range("I113").Value = 106071.25
$ws.Range("J113").Value = 114445
$ws.Range("K113").Value = 106071.25
$ws.Range("L113").Value = 114445
$ws.Range("M113").Value = -103901.25
$ws.Range("N113").Value = -118785
$ws.Range("H132").Value = 3801.0476
$ws.Range("I132").Value = 3443.6924
$ws.Range("J132").Value = 4381.75
$ws.Range("K132").Value = 10331.0772
$ws.Range("L132").Value = 13145.25
$ws.Range("M132").Value = -7801.0772
$ws.Range("N132").Value = -18205.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 11856.429
$ws.Range("I4").Value = 20665
$ws.Range("K4").Value = 20665
$ws.Range("M4").Value = -20552
$ws.Range("H132").Value = 5318.1924
$ws.Range("I132").Value = 5874.4287
$ws.Range("K132").Value = 17623.2861
$ws.Range("M132").Value = -15093.2861
